# Update the "generated at" date/day line.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-11-16 Saturday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2024-11-17 Sunday", 2)

# Update the multiplication-fact table. Addressing cells by (row, column)
# avoids any ambiguity from values that are reused elsewhere in the table
# (e.g. "27×78=2106" is an old value in row 1 but also the new value for
# row 20, column 3).
$table = $d.Tables.Item(1)

$table.Cell(1, 1).Range.Text  = "55×47=2585"
$table.Cell(1, 2).Range.Text  = "37×11=407"
$table.Cell(1, 3).Range.Text  = "86×49=4214"
$table.Cell(1, 4).Range.Text  = "81×38=3078"
$table.Cell(1, 5).Range.Text  = "70×20=1400"

$table.Cell(5, 1).Range.Text  = "54×12=648"
$table.Cell(5, 2).Range.Text  = "36×41=1476"
$table.Cell(5, 3).Range.Text  = "11×52=572"
$table.Cell(5, 4).Range.Text  = "89×82=7298"
$table.Cell(5, 5).Range.Text  = "60×68=4080"

$table.Cell(10, 1).Range.Text = "13×19=247"
$table.Cell(10, 2).Range.Text = "87×86=7482"
$table.Cell(10, 3).Range.Text = "48×99=4752"
$table.Cell(10, 4).Range.Text = "44×68=2992"
$table.Cell(10, 5).Range.Text = "69×36=2484"

$table.Cell(15, 1).Range.Text = "78×71=5538"
$table.Cell(15, 2).Range.Text = "71×29=2059"
$table.Cell(15, 3).Range.Text = "23×50=1150"
$table.Cell(15, 4).Range.Text = "32×26=832"
$table.Cell(15, 5).Range.Text = "58×79=4582"

$table.Cell(20, 1).Range.Text = "49×43=2107"
$table.Cell(20, 2).Range.Text = "60×48=2880"
$table.Cell(20, 3).Range.Text = "27×78=2106"
$table.Cell(20, 4).Range.Text = "77×52=4004"
$table.Cell(20, 5).Range.Text = "22×23=506"
